$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# Update A8: FIRECLASS 64-2 -> FC64-2
$ws.Range("A8").Value = "FC64-2"

# Update A10: FIRECLASS 240-2 -> FC240-2
$ws.Range("A10").Value = "FC240-2"

# Update F8/G8, F9/G9, F10/G10 values - entered as text with leading apostrophe
$ws.Range("F8").Value = "'3.100"
$ws.Range("G8").Value = "'5.000"
$ws.Range("F9").Value = "'3.000"
$ws.Range("G9").Value = "'5.000"
$ws.Range("F10").Value = "'3.100"
$ws.Range("G10").Value = "'5.000"

# Update NGC-1928 value in B4
$ws.Range("B4").Value = "NGC-1928/T957"

# Update selection
$ws.Range("B9").Select()
